$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) -- first worksheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1103
$wsExhibit.Range("F3").Value = 417
$wsExhibit.Range("F4").Value = 1510
$wsExhibit.Range("F5").Value = 8778
$wsExhibit.Range("F6").Value = 96
$wsExhibit.Range("F7").Value = 493
$wsExhibit.Range("F9").Value = 291
$wsExhibit.Range("F12").Value = 11
$wsExhibit.Range("F13").Value = 3627
$wsExhibit.Range("F17").Value = 2187
$wsExhibit.Range("F21").Value = 209
$wsExhibit.Range("F22").Value = 2415

# Sheet "全部类型" (All types) -- fourth worksheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1103
$wsAll.Range("F3").Value = 417
$wsAll.Range("F4").Value = 1510
$wsAll.Range("F5").Value = 8778
$wsAll.Range("F6").Value = 96
$wsAll.Range("F7").Value = 493
$wsAll.Range("F9").Value = 291
$wsAll.Range("F12").Value = 11
$wsAll.Range("F13").Value = 3627
$wsAll.Range("F17").Value = 2188
$wsAll.Range("F21").Value = 209
$wsAll.Range("F22").Value = 2415

$wb.Save()
